$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.615195631980896
$ws.Range("B1").Value = 2.799104928970337
$ws.Range("C1").Value = 3.288551092147827
$ws.Range("D1").Value = 3.676680088043213
$ws.Range("E1").Value = 1.234872102737427
